$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 82: B82 currently holds the text "3" (stored as a string). Convert it
# to a true numeric value of 3, leaving the other cells in the row untouched.
$ws.Range("B82").Value = 3

# New row 83: append a new annotation row for Ying Tang with the given data.
$ws.Range("A83").Value = "Ying Tang"

# B83 should be stored as TEXT "3" (like the other politeness_score text
# cells), not as a number, so force text formatting before assigning it,
# then reset the style so no stray number-format style sticks to the cell.
$ws.Range("B83").NumberFormat = "@"
$ws.Range("B83").Value = "3"
$ws.Range("B83").Style = "Normal"

$ws.Range("C83").Value = "无"
$ws.Range("D83").Value = "FBK"
$ws.Range("E83").Value = "EXP"
$ws.Range("F83").Value = "21c11312-d736-4194-815f-bf7208ef5d55"
$ws.Range("G83").Value = "SJ60SbW0b_annotated.xlsx"
$ws.Range("H83").Value = "These sections include new experiments that illustrate the effect of varying the beta hyperparameter, demonstrate the strength of our approach on the larger scale Inception network for the ILSVRC 2014 classification challenge, and further highlight the effectiveness of our approach in diagnosing model failure modes."
